$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1679.1628
$ws.Range("I15").Value = 1679.1628
$ws.Range("K15").Value = 5037.4884
$ws.Range("M15").Value = -4868.4884

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5374.1665
$ws.Range("I64").Value = 4950
$ws.Range("J64").Value = 5586.25
$ws.Range("K64").Value = 4950
$ws.Range("L64").Value = 5586.25
$ws.Range("M64").Value = -4702
$ws.Range("N64").Value = -6082.25

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5374.1665
$ws.Range("I67").Value = 4950
$ws.Range("J67").Value = 5586.25
$ws.Range("K67").Value = 4950
$ws.Range("L67").Value = 5586.25
$ws.Range("M67").Value = -4092
$ws.Range("N67").Value = -7302.25

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2528.6
$ws.Range("J70").Value = 2531.7778
$ws.Range("L70").Value = 7595.3334
$ws.Range("N70").Value = -8135.3334

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2528.6
$ws.Range("J73").Value = 2531.7778
$ws.Range("L73").Value = 7595.3334
$ws.Range("N73").Value = -9467.3334

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 849.3333
$ws.Range("I80").Value = 705.5
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 2116.5
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -1118.5
$ws.Range("N80").Value = -7996

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 849.3333
$ws.Range("I83").Value = 705.5
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 6349.5
$ws.Range("L83").Value = 18000
$ws.Range("M83").Value = -1357.5
$ws.Range("N83").Value = -27984

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1665.9166
$ws.Range("I107").Value = 1609.6666
$ws.Range("K107").Value = 1609.6666
$ws.Range("M107").Value = 310.3334

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2310.6
$ws.Range("J112").Value = 2389.7778
$ws.Range("L112").Value = 7169.3334
$ws.Range("N112").Value = -9385.3334

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3924.611
$ws.Range("I116").Value = 3795.182
$ws.Range("J116").Value = 4128
$ws.Range("K116").Value = 3795.182
$ws.Range("L116").Value = 4128
$ws.Range("M116").Value = -353.1819999999998
$ws.Range("N116").Value = -11012

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4221.2764
$ws.Range("I132").Value = 3520.525
$ws.Range("K132").Value = 10561.575
$ws.Range("M132").Value = -8031.575000000001

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 698.1667
$ws.Range("I135").Value = 677
$ws.Range("J135").Value = 899.25
$ws.Range("K135").Value = 6093
$ws.Range("L135").Value = 8093.25
$ws.Range("M135").Value = -3558
$ws.Range("N135").Value = -13163.25

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3160.2693
$ws.Range("I137").Value = 2755.8572
$ws.Range("J137").Value = 3632.0833
$ws.Range("K137").Value = 8267.571599999999
$ws.Range("L137").Value = 10896.2499
$ws.Range("M137").Value = -5717.571599999999
$ws.Range("N137").Value = -15996.2499

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1105.2826
$ws.Range("I138").Value = 853
$ws.Range("J138").Value = 2510.8572
$ws.Range("K138").Value = 2559
$ws.Range("L138").Value = 7532.571599999999
$ws.Range("M138").Value = 2581
$ws.Range("N138").Value = -17812.5716

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32416.574
$ws.Range("I32").Value = 33742.145
$ws.Range("K32").Value = 33742.145
$ws.Range("M32").Value = -33455.145

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7414.85
$ws.Range("J45").Value = 3934.6
$ws.Range("L45").Value = 3934.6
$ws.Range("N45").Value = -4688.6

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1744.1333
$ws.Range("I102").Value = 1560.7273
$ws.Range("K102").Value = 1560.7273
$ws.Range("M102").Value = 61.27269999999999

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2748.2942
$ws.Range("I110").Value = 2448.1333
$ws.Range("K110").Value = 2448.1333
$ws.Range("M110").Value = -403.1333

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 45246.74
$ws.Range("I132").Value = 45246.74
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 135740.22
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -133210.22
$ws.Range("N132").Value = $null

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2994.0435
$ws.Range("I86").Value = 2433.0667
$ws.Range("J86").Value = 4045.875
$ws.Range("K86").Value = 2433.0667
$ws.Range("L86").Value = 4045.875
$ws.Range("M86").Value = -1310.0667
$ws.Range("N86").Value = -6291.875

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2994.0435
$ws.Range("I89").Value = 2433.0667
$ws.Range("J89").Value = 4045.875
$ws.Range("K89").Value = 12165.3335
$ws.Range("L89").Value = 20229.375
$ws.Range("M89").Value = -6549.333499999999
$ws.Range("N89").Value = -31461.375

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1783.92
$ws.Range("I94").Value = 1217
$ws.Range("K94").Value = 1217
$ws.Range("M94").Value = -766

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3957.1052
$ws.Range("I107").Value = 3729.923
$ws.Range("J107").Value = 4449.3335
$ws.Range("K107").Value = 3729.923
$ws.Range("L107").Value = 4449.3335
$ws.Range("M107").Value = -1809.923
$ws.Range("N107").Value = -8289.333500000001

# CRP row 17
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 8999
$ws.Range("I17").Value = 8999
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 8999
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -8825
$ws.Range("N17").Value = $null

# CRP row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1025
$ws.Range("I25").Value = 1025
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1025
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -851
$ws.Range("N25").Value = $null

# CRP row 48
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 29998
$ws.Range("J48").Value = 29998
$ws.Range("L48").Value = 29998
$ws.Range("N48").Value = -30950

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4180
$ws.Range("I86").Value = 3987
$ws.Range("K86").Value = 3987
$ws.Range("M86").Value = -2864

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 4180
$ws.Range("I89").Value = 3987
$ws.Range("K89").Value = 19935
$ws.Range("M89").Value = -14319

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 23457.111
$ws.Range("I105").Value = 23457.111
$ws.Range("K105").Value = 23457.111
$ws.Range("M105").Value = -21710.111

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 363.20834
$ws.Range("I107").Value = 343.95456
$ws.Range("J107").Value = 575
$ws.Range("K107").Value = 343.95456
$ws.Range("L107").Value = 575
$ws.Range("M107").Value = 1576.04544
$ws.Range("N107").Value = -4415

# CRP row 108
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 74997
$ws.Range("J108").Value = 74997
$ws.Range("L108").Value = 74997
$ws.Range("N108").Value = -82677

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 774.75
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 774.75
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 2324.25
$ws.Range("M38").Value = $null
$ws.Range("N38").Value = -3018.25

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5887233.5
$ws.Range("I131").Value = 1354
$ws.Range("J131").Value = 10007349
$ws.Range("K131").Value = 4062
$ws.Range("L131").Value = 30022047
$ws.Range("M131").Value = 978
$ws.Range("N131").Value = -30032127

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1744.1875
$ws.Range("I97").Value = 1364
$ws.Range("J97").Value = 2377.8333
$ws.Range("K97").Value = 1364
$ws.Range("L97").Value = 2377.8333
$ws.Range("M97").Value = -868
$ws.Range("N97").Value = -3369.8333

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 17229
$ws.Range("I102").Value = 4959.5
$ws.Range("K102").Value = 4959.5
$ws.Range("M102").Value = -3337.5

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 156964.08
$ws.Range("I113").Value = 128041.375
$ws.Range("K113").Value = 128041.375
$ws.Range("M113").Value = -125871.375

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 31632.588
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 658.4
$ws.Range("I22").Value = 658.4
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 658.4
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -363.4
$ws.Range("N22").Value = $null

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 658.4
$ws.Range("I27").Value = 658.4
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 658.4
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -551.4
$ws.Range("N27").Value = $null

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4209.778
$ws.Range("I100").Value = 3877
$ws.Range("J100").Value = 4476
$ws.Range("K100").Value = 3877
$ws.Range("L100").Value = 4476
$ws.Range("M100").Value = -3336
$ws.Range("N100").Value = -5558

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 74246.17999999999
$ws.Range("I132").Value = 88663.28999999999
$ws.Range("J132").Value = 6966.3335
$ws.Range("K132").Value = 265989.87
$ws.Range("L132").Value = 20899.0005
$ws.Range("M132").Value = -263459.87
$ws.Range("N132").Value = -25959.0005

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2078.4443
$ws.Range("I107").Value = 1000.2
$ws.Range("K107").Value = 3000.6
$ws.Range("M107").Value = -1080.6

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 40123.18
$ws.Range("I132").Value = 46064.625
$ws.Range("K132").Value = 138193.875
$ws.Range("M132").Value = -135663.875

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2459.5862
$ws.Range("I136").Value = 2311.4285
$ws.Range("J136").Value = 2848.5
$ws.Range("K136").Value = 6934.2855
$ws.Range("L136").Value = 8545.5
$ws.Range("M136").Value = -4384.2855
$ws.Range("N136").Value = -13645.5
